# Update FamilySituation STU3 mappings for the GGZ usecase.
# The "Observation.subject(Patient) revinclude relatedperson" / RelatedPerson
# based mapping is replaced with an Observation.component based mapping, and
# a note is added clarifying that marital status is captured through the
# Patient resource.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ZIB 2017")

# Row 4 (Gezinssamenstelling / FamilyComposition): add a note in the Notes
# column (Q) explaining marital status is captured through the Patient.
$ws.Range("Q4").Value2 = "Marital status is captured through the Patient"

# Row 11 (Geboortedatum / DateOfBirth): update the "Maps to" mapping and
# clear the old "Notes" about RelatedPerson / Resident since it is no
# longer relevant with the new Observation.component mapping.
$ws.Range("P11").Value2 = "Observation.component.extension"

# Row 9 (Zorgtaak / CareResponsibility): update the "Maps to" mapping.
$ws.Range("P9").Value2 = "Observation.component"

# Row 10 (Kind / Child): update the "Maps to" mapping.
$ws.Range("P10").Value2 = "Observation.component.value[x]"

$ws.Range("Q11").ClearContents()

# Restore the last active selection to match the author's saved view state.
$ws.Activate()
$ws.Range("N20").Select()
